# [LPF-879]: CCMS Third party report
# Remove the "By Source and Expenditure type", "Provider Contigency" and "MAIN"
# worksheets (and the now-unused "SourceNType" pivot table that lived on the
# "By Source and Expenditure type" sheet), leaving only "Summary" and
# "Transparency Rec".

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Explicitly drop the pivot table that only exists on the sheet we are about
# to remove, so no orphaned pivot table definition is left behind in the
# workbook once the sheet itself is deleted.
$sourceSheet = $wb.Worksheets.Item("By Source and Expenditure type")
foreach ($pt in $sourceSheet.PivotTables()) {
    $pt.TableRange2.Delete() | Out-Null
}

$sourceSheet.Delete() | Out-Null
$wb.Worksheets.Item("Provider Contigency").Delete() | Out-Null
$wb.Worksheets.Item("MAIN").Delete() | Out-Null
